# Add the next three days of vaccination-county-URL data to the log sheet.
# (Mirrors the author's daily "Add files via upload" append of 2021-04-23..25.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New daily rows (29, 30, 31) --------------------------------------
# Column A continues the running "=<prev>+1" date series that fills A2:A28.
# Copy the last existing row's formatting down first (keeps the same date
# style / cellXf index instead of minting a new one), then fill in the
# per-row formula and the new URL for column B.

$ws.Range("A28").Copy()
$ws.Range("A29:A31").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A29").Formula = "=A28+1"
$ws.Range("A30").Formula = "=A29+1"
$ws.Range("A31").Formula = "=A30+1"

$ws.Range("B29").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/6be51bd4348df57c1533fd1a13d3e0fcdd0107c7/counties.json"
$ws.Range("B30").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/c9a530dc6d085617a9ea6d6669b4f9ef8ba3fd50/counties.json"
$ws.Range("B31").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/518c8623cf6c257adbc938fe3ebcce965d2df854/counties.json"

# --- Cosmetic touch-ups that accompanied the upload --------------------
# Column A was manually widened a bit (and no longer auto "best fit").
$ws.Columns.Item(1).ColumnWidth = 12.8

# Selection / scroll position left by the editor at save time.
$ws.Range("F42").Select() | Out-Null
